$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.774.09"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.626.93"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.88"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5078"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2561"
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06335"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.42"
$ws.Range("E10").Value = "  -2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07775"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.262"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.623.16"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.849.92"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5559"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₅7472"
$ws.Range("E17").Value = "  -3.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.794.52"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.418"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.37"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.774"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.980"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.861"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.81"
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1244"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.729"
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.45"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.235"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04876"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.303"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.177"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.546"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8940"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5507"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.530"
$ws.Range("E38").Value = "  -2.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.121.36"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01547"
$ws.Range("E40").Value = "  -1.71%  "
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.546"
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7930"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("E44").Value = "  -2.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.773.97"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -7.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4418"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.63"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05109"
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.583"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9985"
$ws.Range("E51").Value = "  -0.64%  "
